# Re-upload / refresh of the "grupos" standings sheet.
# - Updates the standings table (Po/GP/GC/Sa/Jo columns) for the first
#   four teams in group A after newly-played matches.
# - Widens column F slightly to fit the new "Sa" (goal-difference) values.
# - Leaves the active selection on H5, matching the author's last click
#   before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Real Carira) ---
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 2

# --- Row 3 (Atlântico) ---
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2

# --- Row 4 (Galo Futsal) ---
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = -11
$ws.Range("G4").Value = 2

# --- Row 5 (Vasquinho) ---
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2

# Column F needs to be a bit wider now that it holds "-11".
$ws.Columns.Item(6).ColumnWidth = 2.8333333333333335

# Restore the cursor to the last-selected cell.
$ws.Range("H5").Select()
